# "Update countries & provincias Spain"
#
# Several new country rows were inserted into the original data set (e.g.
# Croacia, San Cristobal y Nieves, Gabon, Benin, San Martin (Parte
# Holandesa), Fiyi, San Bartolome, Islas Virgenes Britanicas, Anguila,
# Guinea-Bisau, Timor Oriental) which pushed the countries below them down
# by one row without touching their statistics. Also the Covid-19 case
# counters for several countries/provinces were refreshed, and the "last
# updated" footer timestamp moved from 13:50 to 14:20.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (country/region name) text that shifted because of the inserted
# rows, plus the footer timestamp in A1.
$textUpdates = @{
    "A1" = "Datos actualizados a 30 de Marzo de 2020 a las 14:20"
    "A52" = "Croacia"
    "A53" = "Eslovenia"
    "A54" = "Serbia"
    "A55" = "Estonia"
    "A176" = "San Cristobal y Nieves"
    "A177" = "Antigua y Barbuda"
    "A179" = "Zimbabue"
    "A180" = "Angola"
    "A181" = "Benin"
    "A182" = "San Martin (Parte Holandesa)"
    "A183" = "Santa Sede"
    "A185" = "Sudan"
    "A186" = "Fiyi"
    "A187" = "San Bartolome"
    "A188" = "Montserrat"
    "A189" = "Nepal"
    "A190" = "Mauritania"
    "A191" = "Butan"
    "A192" = "Islas Turcas y Caicos"
    "A193" = "Nicaragua"
    "A194" = "Gambia"
    "A195" = "Republica de Africa Central"
    "A196" = "Liberia"
    "A197" = "Somalia"
    "A198" = "Republica del Chad"
    "A199" = "Islas Virgenes Britanicas"
    "A200" = "Anguila"
    "A201" = "Guinea-Bisau"
    "A202" = "Belice"
    "A203" = "Timor Oriental"
    "A204" = "Papua Nueva Guinea"
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Updated statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) for the affected rows.
$numericUpdates = @{
    "B12" = 15475
    "C12" = 646
    "E12" = 13325
    "G12" = 27
    "H12" = 327
    "B14" = 11750
    "C14" = 884
    "E14" = 10636
    "G14" = 93
    "H14" = 864
    "F18" = 164
    "B20" = 4393
    "C20" = 109
    "D20" = 12
    "E20" = 4350
    "D22" = 120
    "E22" = 4057
    "B24" = 4028
    "C24" = 328
    "E24" = 3866
    "F24" = 306
    "G24" = 36
    "H24" = 146
    "B40" = 1352
    "C40" = 112
    "E40" = 1329
    "G40" = 2
    "H40" = 13
    "B52" = 790
    "C52" = 77
    "D52" = 55
    "E52" = 729
    "F52" = 26
    "H52" = 6
    "B53" = 756
    "C53" = 26
    "D53" = 10
    "E53" = 735
    "F53" = 28
    "G53" = 0
    "H53" = 11
    "B54" = 741
    "C54" = 0
    "D54" = 42
    "E54" = 685
    "F54" = 25
    "G54" = 1
    "H54" = 14
    "C55" = 36
    "D55" = 20
    "E55" = 692
    "F55" = 10
    "H55" = 3
    "D64" = 14
    "E64" = 473
    "G64" = 3
    "H64" = 29
    "B86" = 230
    "C86" = 6
    "D86" = 13
    "E86" = 192
    "G86" = 3
    "H86" = 25
    "C176" = 5
    "E177" = 7
    "H177" = 0
    "E179" = 6
    "H179" = 1
    "B180" = 7
    "E180" = 5
    "H180" = 2
    "E183" = 6
    "H183" = 0
    "B185" = 6
    "E185" = 4
    "G185" = 1
    "H185" = 2
    "D188" = 0
    "E188" = 5
    "D189" = 1
    "E189" = 4
    "B190" = 5
    "D190" = 2
    "E190" = 3
    "E192" = 4
    "H192" = 0
    "B194" = 4
    "H194" = 1
    "B198" = 3
    "E198" = 3
}
foreach ($ref in $numericUpdates.Keys) {
    $ws.Range($ref).Value = $numericUpdates[$ref]
}
